# Weekly update: insert a new data row for Espinaca (Feria Lagunitas de
# Puerto Montt) above the current row 42, shifting the existing rows
# 42-50 down to 43-51, and populate the new row with the latest week's
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42; this pushes rows 42..50 down to 43..51
# and carries the number formatting of the row above (date format on D).
$ws.Rows.Item(42).Insert()

$ws.Cells.Item(42, 1).Value  = 4
$ws.Cells.Item(42, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(42, 3).Value  = "Los Lagos"
$ws.Cells.Item(42, 4).Value  = 44841
$ws.Cells.Item(42, 5).Value  = 10
$ws.Cells.Item(42, 6).Value  = 100112012
$ws.Cells.Item(42, 7).Value  = "Espinaca"
$ws.Cells.Item(42, 8).Value  = "Sin especificar"
$ws.Cells.Item(42, 9).Value  = "Primera"
$ws.Cells.Item(42, 10).Value = 35
$ws.Cells.Item(42, 11).Value = 12000
$ws.Cells.Item(42, 12).Value = 12000
$ws.Cells.Item(42, 13).Value = 12000
$ws.Cells.Item(42, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(42, 15).Value = "Región Metropolitana"
$ws.Cells.Item(42, 16).Value = 1200
$ws.Cells.Item(42, 17).Value = 10
$ws.Cells.Item(42, 18).Value = "Hortaliza"
